# Update column G ("K" - strikeouts) values in Sheet1 to regenerated figures
# computed from the refreshed source data (replacing the old Strike# based count).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$gValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 0
    6 = 1
    7 = 2
    8 = 1
    9 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 2
    15 = 2
    16 = 2
    17 = 1
    18 = 0
    19 = 0
    20 = 3
    21 = 0
    22 = 3
    23 = 1
    24 = 1
    25 = 2
    26 = 2
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 0
    32 = 2
    33 = 1
    34 = 0
    35 = 1
    36 = 1
    37 = 2
    38 = 0
    39 = 1
    40 = 0
    41 = 2
    42 = 0
    43 = 1
    44 = 0
    45 = 0
    46 = 2
    47 = 1
    48 = 0
    49 = 1
    50 = 1
    51 = 1
    52 = 2
    53 = 3
    54 = 0
    55 = 2
    56 = 3
    57 = 1
    58 = 5
    59 = 3
    60 = 3
    61 = 2
    62 = 2
    63 = 4
    64 = 0
    65 = 2
    66 = 0
    68 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}
